$wb = $excel.ActiveWorkbook

# Sheets 1-4 share the same structure: a "Fonte/Tecnologia" header in A1,
# accented technology names in A2:A12, and loss of the bold/border style
# on A2:A12 (keeping it only on row 1 and on the year headers in B1:E1).
$sheetNames = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A1").Value = "Fonte/Tecnologia"
    $ws.Range("A1").Style = $ws.Range("B1").Style

    $ws.Range("A2").Value = "Hidro"
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A5").Value = "Nuclear"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A7").Value = "Biomassa"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A9").Value = "Solar"
    $ws.Range("A10").Value = "Outros"
    $ws.Range("A11").Value = "Pot. Compl."
    $ws.Range("A12").Value = "GD"

    $ws.Range("A2:A12").Style = "Normal"
}

# Sheet 5: "Emissoes Totais (MtCO2eq)" gets a "Período" header, accented
# row labels, and loses the "Teto" row entirely.
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
$ws5.Range("A1").Value = "Período"
$ws5.Range("A1").Style = $ws5.Range("B1").Style
$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"
$ws5.Range("A2:A3").Style = "Normal"
$ws5.Rows.Item(4).Delete()

# Sheet 6: "Custo Total (bilhões de R$)" gets a "Tipo Expansão" header,
# the cost column header becomes "2015", labels get accents, and the
# values are updated.
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Range("A1").Value = "Tipo Expansão"
$ws6.Range("A1").Style = $ws6.Range("B1").Style
$ws6.Range("B1").Value = "2015"
$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 593
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
$ws6.Range("A2:A3").Style = "Normal"
